# Updated RAD Test Cases for Filing Year drop down and MD CRN changes.
# The "Date" column (column B) on the RAD test-data sheet holds a per-row
# execution timestamp. This re-run stamps each existing test row (except
# rows 19-24, a separate/earlier batch that is left untouched) with a new
# "Sat Feb 17 2024" timestamp sequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timestamps, in row order: rows 2-18 (17 rows) followed by rows 25-48
# (24 rows) = 41 values total.
$dates = @(
  "Sat Feb 17 22:31:09 EST 2024",
  "Sat Feb 17 22:31:18 EST 2024",
  "Sat Feb 17 22:31:27 EST 2024",
  "Sat Feb 17 22:31:37 EST 2024",
  "Sat Feb 17 22:31:47 EST 2024",
  "Sat Feb 17 22:31:56 EST 2024",
  "Sat Feb 17 22:32:06 EST 2024",
  "Sat Feb 17 22:32:15 EST 2024",
  "Sat Feb 17 22:32:25 EST 2024",
  "Sat Feb 17 22:32:34 EST 2024",
  "Sat Feb 17 22:32:43 EST 2024",
  "Sat Feb 17 22:32:53 EST 2024",
  "Sat Feb 17 22:33:02 EST 2024",
  "Sat Feb 17 22:33:12 EST 2024",
  "Sat Feb 17 22:33:21 EST 2024",
  "Sat Feb 17 22:33:30 EST 2024",
  "Sat Feb 17 22:33:40 EST 2024",
  "Sat Feb 17 22:33:50 EST 2024",
  "Sat Feb 17 22:34:00 EST 2024",
  "Sat Feb 17 22:34:09 EST 2024",
  "Sat Feb 17 22:34:19 EST 2024",
  "Sat Feb 17 22:34:28 EST 2024",
  "Sat Feb 17 22:34:37 EST 2024",
  "Sat Feb 17 22:34:47 EST 2024",
  "Sat Feb 17 22:34:57 EST 2024",
  "Sat Feb 17 22:35:06 EST 2024",
  "Sat Feb 17 22:35:16 EST 2024",
  "Sat Feb 17 22:35:25 EST 2024",
  "Sat Feb 17 22:35:34 EST 2024",
  "Sat Feb 17 22:35:44 EST 2024",
  "Sat Feb 17 22:35:53 EST 2024",
  "Sat Feb 17 22:36:03 EST 2024",
  "Sat Feb 17 22:36:12 EST 2024",
  "Sat Feb 17 22:36:22 EST 2024",
  "Sat Feb 17 22:36:31 EST 2024",
  "Sat Feb 17 22:36:40 EST 2024",
  "Sat Feb 17 22:36:51 EST 2024",
  "Sat Feb 17 22:37:00 EST 2024",
  "Sat Feb 17 22:37:09 EST 2024",
  "Sat Feb 17 22:37:19 EST 2024",
  "Sat Feb 17 22:37:28 EST 2024"
)

# Target rows, in the same order as $dates above.
$rows = @(2..18) + @(25..48)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 2).Value = $dates[$i]
}
